$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
}

$ws.Range('D2').Value = '63.825.44'
$ws.Range('E2').Value = '  +1.47%  '
$ws.Range('D3').Value = '3.431.26'
$ws.Range('E3').Value = '  +2.47%  '
$ws.Range('E4').Value = '  +0.04%  '
Set-TextValue 'D5' '572.51'
$ws.Range('E5').Value = '  +2.86%  '
Set-TextValue 'D6' '156.79'
$ws.Range('E6').Value = '  +2.73%  '
Set-TextValue 'D7' '0.999'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '3.433.90'
$ws.Range('E8').Value = '  +2.39%  '
Set-TextValue 'D9' '0.546'
$ws.Range('E9').Value = '  +2.54%  '
Set-TextValue 'D10' '7.46'
$ws.Range('E10').Value = '  -0.65%  '
$ws.Range('E11').Value = '  +3.81%  '
$ws.Range('E12').Value = '  +0.13%  '
$ws.Range('D13').Value = '4.025.62'
$ws.Range('E13').Value = '  +2.59%  '
$ws.Range('E14').Value = '  -3.13%  '
Set-TextValue 'D15' '0.0000192'
$ws.Range('E15').Value = '  +5.67%  '
Set-TextValue 'D16' '27.27'
$ws.Range('E16').Value = '  +1.31%  '
$ws.Range('D17').Value = '63.969.38'
$ws.Range('E17').Value = '  +1.66%  '
$ws.Range('D18').Value = '3.398.64'
$ws.Range('E18').Value = '  +1.87%  '
$ws.Range('E19').Value = '  -1.96%  '
Set-TextValue 'D20' '14.26'
$ws.Range('E20').Value = '  +3.68%  '
Set-TextValue 'D21' '388.58'
$ws.Range('E21').Value = '  +0.04%  '
Set-TextValue 'D22' '8.28'
$ws.Range('E22').Value = '  -1.95%  '
Set-TextValue 'D23' '1.00'
$ws.Range('E23').Value = '  +0.18%  '
Set-TextValue 'D24' '0.539'
$ws.Range('E24').Value = '  -0.06%  '
Set-TextValue 'D25' '72.37'
$ws.Range('E25').Value = '  +2.48%  '
Set-TextValue 'D26' '0.0000121'
$ws.Range('E26').Value = '  +24.46%  '
Set-TextValue 'D27' '9.58'
$ws.Range('E27').Value = '  +8.84%  '
Set-TextValue 'D28' '0.178'
$ws.Range('E28').Value = '  -1.31%  '
$ws.Range('E29').Value = '  -0.17%  '
Set-TextValue 'D30' '6.09'
$ws.Range('E30').Value = '  +9.05%  '
$ws.Range('E31').Value = '  +5.53%  '
$ws.Range('E32').Value = '  +1.47%  '
$ws.Range('B33').Value = 'RenderToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D33' '6.48'
$ws.Range('E33').Value = '  +1.13%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D34' '23.43'
$ws.Range('E34').Value = '  +1.95%  '
$ws.Range('E35').Value = '  +0.05%  '
Set-TextValue 'D36' '6.93'
$ws.Range('E36').Value = '  +3.12%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D37' '1.47'
$ws.Range('E37').Value = '  -0.62%  '
$ws.Range('B38').Value = 'Monero'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D38' '159.27'
$ws.Range('E38').Value = '  -0.21%  '
$ws.Range('E39').Value = '  +3.46%  '
$ws.Range('D40').Value = '2.925.53'
$ws.Range('E40').Value = '  +3.73%  '
$ws.Range('E41').Value = '  -2.49%  '
$ws.Range('E42').Value = '  -0.53%  '
Set-TextValue 'D43' '0.0318'
$ws.Range('E43').Value = '  +2.27%  '
Set-TextValue 'D44' '4.43'
$ws.Range('E44').Value = '  +2.48%  '
$ws.Range('B45').Value = 'Mantle'
$ws.Range('C45').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D45' '0.767'
$ws.Range('E45').Value = '  +2.64%  '
$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D46' '41.56'
Set-TextValue 'D47' '23.71'
$ws.Range('E47').Value = '  +7.93%  '
Set-TextValue 'D48' '1.09'
$ws.Range('E48').Value = '  +4.36%  '
Set-TextValue 'D49' '2.21'
$ws.Range('E49').Value = '  +21.66%  '
Set-TextValue 'D50' '6.49'
$ws.Range('E50').Value = '  +3.36%  '
Set-TextValue 'D51' '0.848'
$ws.Range('E51').Value = '  +5.01%  '
